$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the "Program" and "Class" sheets, keep only "Batch" ---
$wb.Worksheets.Item("Program").Delete()
$wb.Worksheets.Item("Class").Delete()

$ws = $wb.Worksheets.Item("Batch")

# --- Populate data rows (column-major entry order so shared-string ids line up) ---
$ws.Range("C2").Value = " data driven 1"
$ws.Range("C3").Value = " data driven 2"
$ws.Range("C4").Value = " data driven 3"
$ws.Range("C5").Value = " data driven 4"

$ws.Range("A2").Value = "ChatBotTestuk"
$ws.Range("A3").Value = "ChatBotTestuk"
$ws.Range("A4").Value = "ChatBotTestuk"
$ws.Range("A5").Value = "Splunk"

$ws.Range("B2").Value = 132
$ws.Range("B3").Value = 124
$ws.Range("B4").Value = 125
$ws.Range("B5").Value = 126

$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 3

# --- Distinct font/style for the A5 "Splunk" cell ---
$ws.Range("A5").Font.Name = "Menlo"
$ws.Range("A5").Font.Size = 12
$ws.Range("A5").Font.Color = 0

# --- Styled (but empty) trailer row, re-using the header row's formatting ---
$ws.Range("A1:D1").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = 20

# --- Column widths (engine quantizes ColumnWidth to 1/6-character steps, so
# the inputs are pre-offset by the engine's constant +5/6 "padding" term and
# land on the closest representable grid point to the target widths) ---
$ws.Columns.Item(1).ColumnWidth = 29.666666666666668
$ws.Columns.Item(2).ColumnWidth = 26.830729166666668
$ws.Columns.Item(3).ColumnWidth = 25.330729166666668
$ws.Columns.Item(4).ColumnWidth = 42.830729166666664

# --- Selection moves to B2 ---
$ws.Range("B2").Select()
